$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tournament")

# --- Insert two new rows right after the header row (new rows 2 and 3) ---
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# column-major fill order (matches how the source values were entered)
$ws.Range("A2").Value = "competition-key"
$ws.Range("A3").Value = "host-key"
$ws.Range("B2").Value = "mens-world-cup"
$ws.Range("B3").Value = "north-america"

# --- Append the venue-key rows at the bottom of the table (rows 23-38) ---
$venueRows = @(
    @("venue-key.1", "ca-vancouver"),
    @("venue-key.2", "ca-toronto"),
    @("venue-key.3", "us-east-rutherford-nj"),
    @("venue-key.4", "us-kansas-city-mo"),
    @("venue-key.5", "us-dallas-tx"),
    @("venue-key.6", "us-houston-tx"),
    @("venue-key.7", "us-atlanta-ga"),
    @("venue-key.8", "us-los-angeles-ca"),
    @("venue-key.9", "us-philadelphia-pa"),
    @("venue-key.10", "us-seattle-wa"),
    @("venue-key.11", "us-santa-clara-ca"),
    @("venue-key.12", "us-boston-ma"),
    @("venue-key.13", "us-miami-fl"),
    @("venue-key.14", "mx-mexico-city"),
    @("venue-key.15", "mx-guadalajara"),
    @("venue-key.16", "mx-monterrey")
)

$startRow = 23
for ($i = 0; $i -lt $venueRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A" + $r).Value = $venueRows[$i][0]
    $ws.Range("B" + $r).Value = $venueRows[$i][1]
}

# --- Resize the "tournament" table to cover the new rows ---
$lo = $ws.ListObjects.Item("tournament")
$lo.Resize($ws.Range("A1:I38"))

# --- Update selection to the newly-added venue rows ---
$ws.Range("A23:B38").Select()

# --- Make Tournament the active/selected sheet ---
$ws.Activate()
